# Ontologia_CROMA.xlsx - "Add files via upload" edit
#
# Content changes:
#  - "Classes" sheet (J6): "tem.mistura some Cor.Digital" -> "tem.rgb some
#    Cor.Digital or tem.rgba some Cor.Digital"
#  - "FatosIn" sheet (V6:V16): "mistura" -> "rgb"
#  - The now-unused shared strings "mistura" / "tem.mistura some Cor.Digital"
#    are dropped automatically by the engine once nothing references them
#    (no cell still needs K2:K5's "tem.red / tem.green / tem.blue / tem.alfa"
#    text, which is unchanged content-wise - only its shared-string index
#    shifts as a side effect of the table shrinking).
#
# View/selection changes:
#  - "Classes": selection moves from J6 to J9
#  - "FatosIn": selection moves from X13 to the multi-cell range V6:V16
#    (active cell V6)

$wb = $excel.ActiveWorkbook

$wsClasses = $wb.Worksheets.Item("Classes")
$wsFatosIn = $wb.Worksheets.Item("FatosIn")

# --- FatosIn: "mistura" -> "rgb" for the whole V6:V16 column ---
$wsFatosIn.Activate()
$wsFatosIn.Range("V6:V16").Value = "rgb"
$wsFatosIn.Range("V6:V16").Select()

# --- Classes: J6 text update ---
$wsClasses.Activate()
$wsClasses.Range("J6").Value = "tem.rgb some Cor.Digital or tem.rgba some Cor.Digital"
$wsClasses.Range("J9").Select()
